# Generate Report for Handback
# Update the "Latest Handback DateTime" (column K, row 2) for the
# 7e02af04-... file entry on both the zh-cn and de-de report sheets.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("K2").Value = "2016-11-03 20:09:32"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("K2").Value = "2016-11-03 20:09:49"
